$wb = $excel.ActiveWorkbook

# ---- "index" sheet ----------------------------------------------------
$wsIndex = $wb.Worksheets.Item("index")

# Update the "back to manual index" link to point at the new guide41 top page
$wsIndex.Range("B5").Value = '<a class="btn btn-primary btn-xs" role="button" href="https://support.vle.hiroshima-u.ac.jp/mdl:guide41:top" style="width:45%">マニュアルの目次へ戻る</a><br>'

# Move the remembered selection on this sheet from B6 to B5
[void]$wsIndex.Range("B5").Select()

# ---- "p1" sheet ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("p1")

# Header text now flags the dashboard page as still in preparation
$ws1.Range("B2").Value = "ダッシュボードの構成【準備中】"

# Replace the old "please edit" placeholder with the "under construction" notice
$ws1.Range("B7").Value = "現在作成中です。
順次公開いたしますので、公開前のページについては過去のマニュアルを参照してください。
■過去のマニュアル
https://support.vle.hiroshima-u.ac.jp/files/public/hirodai-moodle-faculty-document-20230306.pdf"
$ws1.Range("B7").RowHeight = 54

# This sheet becomes the active / tab-selected sheet with B2 selected
[void]$ws1.Activate()
[void]$ws1.Range("B2").Select()

# ---- "p2" sheet ---------------------------------------------------------
# No textual changes here; activating "p1" above already removes the
# tabSelected flag that used to live on this sheet, and its own selection
# (B15) is left untouched.
